$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Extend the merged footnote range from A17:G17 to A17:J17 (before touching row 17 content) ---
$ws.Range("A17:G17").UnMerge() | Out-Null
$ws.Range("A17:J17").Merge() | Out-Null

# --- 2) Add columns H, I, J: same width as the existing data columns, and per-row formatting copied from column G ---
$ws.Range("H1:J17").ColumnWidth = $ws.Range("A1").ColumnWidth

for ($r = 1; $r -le 16; $r++) {
    $ws.Range("G$r").Copy() | Out-Null
    $ws.Range("H$r`:J$r").PasteSpecial(-4122) | Out-Null
}
$ws.Range("A16").Copy() | Out-Null
$ws.Range("A17:J17").PasteSpecial(-4122) | Out-Null

# --- 3) Update existing cells (columns A-G) whose values changed ---
$ws.Range("C1").Value = 'FD-FE'
$ws.Range("E1").Value = 'FD ES'
$ws.Range("F1").Value = 'FD-FE ES'
$ws.Range("G1").Value = 'FE ES'
$ws.Range("C2").Value = '-0.0038*'
$ws.Range("B3").Value = '(0.0045)'
$ws.Range("B4").Value = -0.5597
$ws.Range("F4").Value = '-0.8506***'
$ws.Range("B5").Value = '(0.3580)'
$ws.Range("B6").Value = 0.1038
$ws.Range("E6").Value = '0.1185**'
$ws.Range("F6").Value = '-0.2473*'
$ws.Range("G6").Value = '-0.4732+'
$ws.Range("B7").Value = '(0.2566)'
$ws.Range("F7").Value = '(0.1252)'
$ws.Range("G7").Value = '(0.2637)'
$ws.Range("B8").Value = '0.4707*'
$ws.Range("C8").Value = '0.8345*'
$ws.Range("D8").Value = '1.4843*'
$ws.Range("F8").Value = '0.2528+'
$ws.Range("B9").Value = '(0.2376)'
$ws.Range("F9").Value = '(0.1407)'
$ws.Range("G9").Value = '(0.1830)'
$ws.Range("B10").Value = '6.2357***'
$ws.Range("D10").Value = '0.8435*'
$ws.Range("F10").Value = '0.8083***'
$ws.Range("B11").Value = '(0.5858)'
$ws.Range("F11").Value = '(0.0460)'
$ws.Range("B12").Value = '1.6517***'
$ws.Range("B13").Value = '(0.3458)'
$ws.Range("G13").Value = '(0.0712)'
$ws.Range("B14").Value = '0.8143**'
$ws.Range("D14").Value = '-0.5303+'
$ws.Range("B15").Value = '(0.2503)'
$ws.Range("F15").Value = '(0.0350)'
$ws.Range("G15").Value = '(0.0919)'
$ws.Range("B16").Value = 0.539
$ws.Range("A17").Value = '+ p < 0.1, * p < 0.05, ** p < 0.01, *** p < 0.001'

# --- 4) Fill new columns H, I, J (column-major order, matching the source generation order) ---
$ws.Range("H1").Value = 'FD E'
$ws.Range("H2").Value = '0.0310***'
$ws.Range("H3").Value = '(0.0034)'
$ws.Range("H4").Value = '-0.5498***'
$ws.Range("H5").Value = '(0.0430)'
$ws.Range("H6").Value = '0.0835*'
$ws.Range("H7").Value = '(0.0377)'
$ws.Range("H8").Value = '0.0639**'
$ws.Range("H9").Value = '(0.0195)'
$ws.Range("H10").Value = '0.9322***'
$ws.Range("H11").Value = '(0.0371)'
$ws.Range("H12").Value = '0.0690+'
$ws.Range("H13").Value = '(0.0409)'
$ws.Range("H14").Value = -0.0052
$ws.Range("H15").Value = '(0.0199)'
$ws.Range("H16").Value = 0.936

$ws.Range("I1").Value = 'FD-FE E'
$ws.Range("I2").Value = 0.0007
$ws.Range("I3").Value = '(0.0008)'
$ws.Range("I4").Value = '-0.4258***'
$ws.Range("I5").Value = '(0.0408)'
$ws.Range("I6").Value = '-0.0626+'
$ws.Range("I7").Value = '(0.0362)'
$ws.Range("I8").Value = -0.0153
$ws.Range("I9").Value = '(0.0166)'
$ws.Range("I10").Value = '0.8903***'
$ws.Range("I11").Value = '(0.0235)'
$ws.Range("I12").Value = 0.0109
$ws.Range("I13").Value = '(0.0136)'
$ws.Range("I14").Value = 0.0166
$ws.Range("I15").Value = '(0.0121)'
$ws.Range("I16").Value = 0.944

$ws.Range("J1").Value = 'FE E'
$ws.Range("J2").Value = -0.0004
$ws.Range("J3").Value = '(0.0012)'
$ws.Range("J4").Value = '-0.4078***'
$ws.Range("J5").Value = '(0.0832)'
$ws.Range("J6").Value = -0.0972
$ws.Range("J7").Value = '(0.0654)'
$ws.Range("J8").Value = '0.0581+'
$ws.Range("J9").Value = '(0.0320)'
$ws.Range("J10").Value = '0.9058***'
$ws.Range("J11").Value = '(0.0331)'
$ws.Range("J12").Value = '0.0729*'
$ws.Range("J13").Value = '(0.0336)'
$ws.Range("J14").Value = -0.0416
$ws.Range("J15").Value = '(0.0259)'
$ws.Range("J16").Value = 0.961
